$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParagraphStartingWith($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $cand = $d.Paragraphs.Item($i)
        if ($cand.Range.Text -like "$prefix*") {
            return $cand
        }
    }
    return $null
}

# 1) "UI [bloops] for selecting menu options." - "bloops" is wrapped in
#    proofErr spell-check markers across 3 separate runs; collapse them
#    back into a single clean run (no proofErr markup) in that bullet.
$pBloops = Find-ParagraphStartingWith "UI "
$bloopsXml = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">UI bloops for selecting menu options. </w:t></w:r></w:p>'
$pBloops.Range.InsertXML($bloopsXml) | Out-Null

# 2) Remove the "Add a ammo pick up." bullet entirely (text + its paragraph
#    mark), while keeping the bookmark (_GoBack) that lives in that same
#    paragraph, and merging the following "Handle player death more
#    appropriately." bullet's text into it.
$pAmmo = Find-ParagraphStartingWith "Add a ammo pick up."
$pDeath = $pAmmo.Next()
$rng = $d.Range($pAmmo.Range.Start, $pDeath.Range.End - 1)
$rng.Text = "Handle player death more appropriately. "
# The old "Handle player death..." paragraph is now an empty leftover
# paragraph immediately following; delete it (and its mark) to fully
# collapse the two paragraphs into one.
$pAmmo.Next().Range.Delete()

# 3) "Font: Grunge Tank by [NalGames]" - "NalGames" is wrapped in proofErr
#    spell-check markers across 3 separate runs; collapse back into one
#    clean run (no proofErr markup).
$pFont = Find-ParagraphStartingWith "Font: Grunge Tank by"
$fontXml = '<w:p ' + $wNs + '><w:r><w:t>Font: Grunge Tank by NalGames</w:t></w:r></w:p>'
$pFont.Range.InsertXML($fontXml) | Out-Null
